$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A1 text from "Category" to "categories"
$ws.Range("A1").Value = "categories"

# Update selection to A2
$ws.Range("A2").Select()
